$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

For ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldVal) {
        $cell.Value = $newVal
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed 'Recorded By' cells in column G"
